$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Enter the value as a formula producing the literal string, then
    # convert the formula to a static value in place. This keeps the
    # cell a plain text/shared-string value (not auto-coerced to a
    # number/date) while leaving the cell's existing style untouched.
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163) | Out-Null
}

# Row 5: Créditos-aula: 4 -> 2
Set-TextValue $ws.Range("B5") "2"
Set-TextValue $ws.Range("C5") "2"

# Row 7: Carga horária: 60 h -> 30 h
$ws.Range("B7").Value = "30 h"
$ws.Range("C7").Value = "30 h"

# Row 8: Ativação: 01/01/2016 -> 01/01/2023
Set-TextValue $ws.Range("B8") "01/01/2023"
Set-TextValue $ws.Range("C8") "01/01/2023"

# Row 11: Objectives: add the English objectives text (new cells -> inherit
# column B/C style from the adjacent rows instead of the row's default).
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Value = "To introduce new students to an understanding of what a career is and the conceptual bases of Physical Engineering."
$ws.Range("C11").Value = "To introduce new students to an understanding of what a career is and the conceptual bases of Physical Engineering."

# Row 13: Programa resumido: "Semestral" -> "01/01/2023" (matches source data)
Set-TextValue $ws.Range("B13") "01/01/2023"
Set-TextValue $ws.Range("C13") "01/01/2023"

# Row 14: Short syllabus: add text (new cells)
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122) | Out-Null
$ws.Range("B14").Value = "The Physics Engineering career. Basic engineering concepts. Skills and Abilities of an Engineer. Conceptual physics. Realization of experiments and projects of Physical Engineering."
$ws.Range("C14").Value = "The Physics Engineering career. Basic engineering concepts. Skills and Abilities of an Engineer. Conceptual physics. Realization of experiments and projects of Physical Engineering."

# Row 15: Programa: "01/01/2016" -> "519033 - Carlos Yujiro Shigue"
$ws.Range("B15").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C15").Value = "519033 - Carlos Yujiro Shigue"

# Row 16: Syllabus: add text (new cells)
$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("B16").Value = "The career of Engineering Physics. Scientists x engineers: the interdisciplinary role of Engineering Physics. Fields of action.Physics as a conceptual science: How to learn Physics. Realization of demonstrations and significant scientific experiments in Physics.Basic engineering concepts. Skills and competences of an engineer.Development of a thematic project of Physical Engineering.Competition between projects from different groups.Evaluation of competitions and the discipline as a whole."
$ws.Range("C16").Value = "The career of Engineering Physics. Scientists x engineers: the interdisciplinary role of Engineering Physics. Fields of action.Physics as a conceptual science: How to learn Physics. Realization of demonstrations and significant scientific experiments in Physics.Basic engineering concepts. Skills and competences of an engineer.Development of a thematic project of Physical Engineering.Competition between projects from different groups.Evaluation of competitions and the discipline as a whole."

# Row 18: Método: "519033 - Carlos Yujiro Shigue" -> "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("B18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"

$wb.Save()
